$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "Save" column, matching style of the other header cells (B1:G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Data values for the new "Save" column
$values = @(0, 1, 0, 0, 0, 0, 0, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
